# Odias labels workbook update:
#  - insert a new worksheet "19AUG24_LOWAGGLOM_TANDEM" as the first sheet
#  - populate it with a Date / Agglom / da (nm) table
#  - becomes the active / selected sheet (as it was the last one edited)

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add()
$ws.Move($firstSheet)
$ws.Name = "19AUG24_LOWAGGLOM_TANDEM"

# Header row
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Agglom"
$ws.Range("C1").Value = "da"
$ws.Range("A1:C1").NumberFormat = "@"

# Data rows: Date / Agglom / diameter (nm)
$dates = @(
    "19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24",
    "19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24","19 AUG. 24"
)
$diameters = @(31.5, 37, 43.3, 50.8, 59.5, 69.7, 81.7, 95.8, 112.3, 131.6, 154.3, 180.8, 211.9)

for ($i = 0; $i -lt $diameters.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $dates[$i]
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "low"
    $ws.Range("C$row").Value = $diameters[$i]
    $ws.Range("C$row").NumberFormat = "0.00"
}

# Footer row: blank label cell (formatted) + unit label
$ws.Range("A15").NumberFormat = "@"
$ws.Range("C15").Value = "nm"
$ws.Range("C15").NumberFormat = "@"

# Column A width to match the source layout
$ws.Columns.Item(1).ColumnWidth = 15

# This new sheet is the active / selected one, with H7 as the last selection
$ws.Range("H7").Select() | Out-Null
